$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $ref, $val) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '57.524.36'
$ws.Range('E2').Value = '  -4.26%  '
$ws.Range('D3').Value = '2.920.47'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue $ws 'D5' '547.47'
$ws.Range('E5').Value = '  -3.97%  '
Set-TextValue $ws 'D6' '129.50'
$ws.Range('E6').Value = '  +3.68%  '
$ws.Range('E7').Value = '  -0.06%  '
Set-TextValue $ws 'D8' '0.512'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').Value = '2.912.28'
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('E10').Value = '  -4.36%  '
$ws.Range('E11').Value = '  -6.13%  '
Set-TextValue $ws 'D12' '0.444'
$ws.Range('E12').Value = '  +0.55%  '
Set-TextValue $ws 'D13' '0.0000218'
$ws.Range('E13').Value = '  -1.37%  '
Set-TextValue $ws 'D14' '32.57'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '3.402.09'
$ws.Range('E16').Value = '  -2.45%  '
Set-TextValue $ws 'D17' '6.89'
$ws.Range('E17').Value = '  +5.72%  '
$ws.Range('D18').Value = '2.916.37'
$ws.Range('E18').Value = '  -2.53%  '
$ws.Range('D19').Value = '57.534.00'
$ws.Range('E19').Value = '  -4.26%  '
Set-TextValue $ws 'D20' '415.12'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('E22').Value = '  +2.44%  '
$ws.Range('E23').Value = '  -1.40%  '
Set-TextValue $ws 'D24' '13.11'
$ws.Range('E24').Value = '  +1.65%  '
Set-TextValue $ws 'D25' '79.51'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +0.07%  '
Set-TextValue $ws 'D28' '2.45'
$ws.Range('E28').Value = '  -2.67%  '
$ws.Range('E29').Value = '  +1.85%  '
Set-TextValue $ws 'D30' '7.37'
$ws.Range('E30').Value = '  +2.50%  '
Set-TextValue $ws 'D31' '25.13'
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('E32').Value = '  -2.37%  '
Set-TextValue $ws 'D33' '0.0963'
$ws.Range('E33').Value = '  +1.61%  '
Set-TextValue $ws 'D34' '5.66'
$ws.Range('E34').Value = '  +1.16%  '
Set-TextValue $ws 'D35' '0.927'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  +1.24%  '
Set-TextValue $ws 'D37' '48.15'
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('D39').Value = '0.0₃0674'
$ws.Range('E39').Value = '  +2.70%  '
Set-TextValue $ws 'D40' '2.55'
$ws.Range('E40').Value = '  +4.67%  '
$ws.Range('E41').Value = '  -0.44%  '
Set-TextValue $ws 'D42' '375.52'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.692.72'
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D44' '0.0342'
$ws.Range('E44').Value = '  -3.65%  '
Set-TextValue $ws 'D46' '123.49'
$ws.Range('E46').Value = '  +1.86%  '
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('E49').Value = '  -1.48%  '
Set-TextValue $ws 'D50' '22.83'
$ws.Range('E50').Value = '  -1.82%  '
$ws.Range('E51').Value = '  -0.87%  '
